# Mise à jour de l'attribution des tâches
# Fill in the "qui le fait" (who does it) column (D) for several checklist
# rows, and update two existing entries from a single initial to a
# combined "X/Y" pair of initials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly attributed tasks (column D was empty before)
$ws.Range("D19").Value = "X/S"
$ws.Range("D20").Value = "J"
$ws.Range("D22").Value = "H"
$ws.Range("D42").Value = "X/S"

# Re-attributed tasks (column D previously held a single initial "J")
$ws.Range("D48").Value = "X/J"
$ws.Range("D50").Value = "H/J"

# Row 38 was resized (its wrapped text now needs less vertical space)
$ws.Rows.Item(38).RowHeight = 31.7

# Leave the selection / scroll position where the author ended up
$excel.ActiveWindow.TopLeftCell = $ws.Range("A6")
$ws.Range("D21").Select()
